# Homework exercise update: add new glossary terms and fix a couple of
# leftover "centered+font" styles that are now equivalent to the plain
# centered style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix stray duplicate style on B7 / A28 (now identical to style used
# everywhere else: centered, default font) -----------------------------
$ws.Range("B7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A28").HorizontalAlignment = -4108  # xlCenter

# --- Append the new glossary rows (39-58) ------------------------------
# Row 39 was authored translation-first (B then A), row 40 term-first
# (A then B), so the shared-string table ends up in that exact order.
$ws.Cells.Item(39, 2).Value = "livello"
$ws.Cells.Item(39, 1).Value = "layer"

$ws.Cells.Item(40, 1).Value = "protocol"
$ws.Cells.Item(40, 2).Value = "protocollo"

$terms = @(
    "header",
    "legge di metcalfe",
    "payload",
    "livelli iso/osi",
    "livelli tcp/ip",
    "flow control",
    "data order",
    "ieee vari",
    "spoofing",
    "multicast",
    "ethernet",
    "type field",
    "vlan tag",
    "sfd",
    "bande",
    "BSS",
    "ESS",
    "datagrammi"
)

$row = 41
foreach ($term in $terms) {
    $ws.Cells.Item($row, 1).Value = $term
    $row = $row + 1
}

# --- Keep the view roughly where the author left it --------------------
$ws.Range("A58").Select()
